# Review_430.docx edit: swap in the 30.03.25 "Efficient Online Data Mixing"
# review body, in place of the 01.04.25 "OPTIMIZING PRETRAINING DATA
# MIXTURES" review, and append the new closing paragraphs / link.
$d = $word.ActiveDocument

# --- Paragraph 1: date + title (two runs split by a manual line break) ---
$d.Paragraphs.Item(1).Range.Text = "המאמר היומי של מייק - 30.03.25" + [char]11 + "Efficient Online Data Mixing For Language Model Pre-Training"

# --- Paragraphs 2-6: rewritten body text ---
$d.Paragraphs.Item(2).Range.Text = "ממשיך בלסקור את קו המחקר בנושא אופטימיזציה של אימון מודלים (בפרט מודלי שפה) כאשר יש בידינו כמה דאטהסטים שונים. מכיוון שכבר הגדרתי את הבעיה בסקירות מ-26.03 ומ-28.03 לא אעשה זאת בסקירה זו ומיד אתחיל מהסבר הרעיון העיקרי של המאמר."
$d.Paragraphs.Item(3).Range.Text = "המאמר ניגש לבעיה בצורה שונה משני המאמרים הקודמים שסקרתי אך לדעתי (למרות המורכבות מתמטית מסוימת) הגישה המוצעת במאמר היא מאוד אינטואיטיבית. המחברים מנסים לפתור בעיית בניית דאטהסט D לאימון מודל שפה מהדאטהסטים D_1,..., D_n עם מה שנקרא Multi-Arm Bandits או MAB. אזכיר בעיית MBA מוגדרת באופן הבא: יש לנו כמה n מכונות מזל עם הסתברויות זכייה p1,..., p_n שלא ידועות לנו מראש. המטרה היא למצוא אסטרטגית בחירת מכונה הממקסמת את הזכייה (נגיד, התוחלת שלה) כאשר יש לנו N נסיונות. "
$d.Paragraphs.Item(4).Range.Text = "שימו לב שבעיית אופטימיזצית האימון שלנו די דומה ל-MBA - גם פה אנו צריכים למצוא את אסטרטגית בחירת דאטהסטים לאימון בלי שאנו יודעים מה ״ההשפעה״ של כל דאטהסט לתוצאת האימון הסופית. בלי להיכנס יותר מדי עמוק למתמטיקה (תהליך החלטה מרקובי, התפלגות גיבס וכדומה) המטרה למצוא התפלגות p1,..., p_n על הדאטהסטים שלנו למקסום ביצועי המודל המאומן. הקאץ' כאן שהתפלגות זו משתנה עם האיטרציות כאשר איטרציה במקרה הזה היא צעד אחד (או מספר כלשהו אך קבוע מראש) על הדאטה מהדאטהסט D_i שנבחר באיטרציה זו. "
$d.Paragraphs.Item(5).Range.Text = "כלומר כל פעם אנו בוחרים דאטהסט עם ההתפלגות הנוכחית p, מאמנים את המודל על הדאטה מהדאטהסט הנבחר ומעדכנים את p בהתבסס על תוצאות האימון. כמובן נשאלת השאלה איך ניתן לקבוע p עבור איטרציה הבא על סמך התוצאות של האיטרציה(בחירת דאטהסט) הקודמת. וכאן אנו מגיעים למה שנקרא תגמול (reward) שהוא משקף את ״ההצלחה״ בבחירת הדאטהסט d_i באיטרציה זו. אם האימון על d_i היה מוצלח, אנו רוצה להגדיל את ההסתברות שלו (על חשבון האחרים) כאשר אם הוא פחות מוצלח אז צריך להקטין אותה. "
$d.Paragraphs.Item(6).Range.Text = "אוקיי, אז מה זה בעצם התגמול כאן? התגמול כאן היא המידה שהמודל ירוויח מהדאטה מדאטהסט D_i כלומר ילמד יותר סוג של information gain או IG. המאמר מחשב את IG בתור פרפלקסיטי (שזה אקספוננט של הלוס) על הדאטה של דאטהסט d_i. לוס זה משוערך על סמך באץ' מהדאטהסט. בנוסף יש גם עניין של exploration כי אנו לא רוצים ״להקטין דרסטית״ את הסתברות בחירה של דאטהסט מסוים על סמך מעט באצ'ים ואז מגדילים (כמו שמקובל ב-MBA ובשיטות אחרות של RL) כל הסתברות p_i במספר קטן ε_t שיורד עם האיטרציות."

# --- Paragraph 7: the old arxiv link becomes the "3 steps" lead-in ---
$d.Paragraphs.Item(7).Range.Text = "אז האלגוריתם הסופי מכיל 3 שלבים:"

# --- New paragraphs appended after paragraph 7 (the 3-step list, outro, new link) ---
$d.Paragraphs.Item(7).Range.InsertParagraphAfter()
$d.Paragraphs.Item(8).Range.Text = "עדכון הסתברויות בחירה p1,..., p_n"

$d.Paragraphs.Item(8).Range.InsertParagraphAfter()
$d.Paragraphs.Item(9).Range.Text = "דגימת דאטה מהדאטהסטים  D_1,..., D_n לפי הסתברויות אלה ואימון מודל על דאטה"

$d.Paragraphs.Item(9).Range.InsertParagraphAfter()
$d.Paragraphs.Item(10).Range.Text = "עדכון נוסף של ההסתברויות בהתבסס על המודל המאומן בשלב 2" + [char]11

$d.Paragraphs.Item(10).Range.InsertParagraphAfter()
$d.Paragraphs.Item(11).Range.Text = "מאמר מומלץ - נהניתי לצלול אליו "

$d.Paragraphs.Item(11).Range.InsertParagraphAfter()
$d.Paragraphs.Item(12).Range.Text = "https://arxiv.org/pdf/2312.02406" + [char]11

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
